$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 82, pushing existing rows 82-86 down to 83-87
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row 82 with the new weekly price observation
$ws.Cells.Item(82, 1).Value = 4
$ws.Cells.Item(82, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(82, 3).Value = "Los Lagos"
$ws.Cells.Item(82, 4).Value = 44747
$ws.Cells.Item(82, 5).Value = 10
$ws.Cells.Item(82, 6).Value = 100112031
$ws.Cells.Item(82, 7).Value = "Poroto verde"
$ws.Cells.Item(82, 8).Value = "Magnum"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 35
$ws.Cells.Item(82, 11).Value = 27000
$ws.Cells.Item(82, 12).Value = 27000
$ws.Cells.Item(82, 13).Value = 27000
$ws.Cells.Item(82, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(82, 15).Value = "Perú"
$ws.Cells.Item(82, 16).Value = 1080
$ws.Cells.Item(82, 17).Value = 25
$ws.Cells.Item(82, 18).Value = "Hortaliza"
